$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$win = $excel.ActiveWindow
$props = $win | Get-Member
Write-Host $props
